$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.974.91"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "1.646.03"
$ws.Range("E3").Value = "  -1.30%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.37"
$ws.Range("E5").Value = "  +2.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5229"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2608"
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06367"
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.70"
$ws.Range("E10").Value = "  -1.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07694"
$ws.Range("E11").Value = "  +2.00%  "
$ws.Range("D12").Value = "1.641.67"
$ws.Range("E12").Value = "  -1.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.427"
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("D14").Value = "1.868.61"
$ws.Range("E14").Value = "  -1.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5523"
$ws.Range("E15").Value = "  +1.86%  "
$ws.Range("D16").Value = "0.0₅8249"
$ws.Range("E16").Value = "  +3.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.71"
$ws.Range("E17").Value = "  -2.62%  "
$ws.Range("D18").Value = "25.987.15"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.706"
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "189.44"
$ws.Range("E21").Value = "  +1.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.16"
$ws.Range("E22").Value = "  -0.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.263"
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.01"
$ws.Range("E25").Value = "  -3.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1253"
$ws.Range("E26").Value = "  +1.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.390"
$ws.Range("E27").Value = "  -0.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.91"
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.393"
$ws.Range("E29").Value = "  +2.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05895"
$ws.Range("E30").Value = "  -6.14%  "
$ws.Range("E31").Value = "  -0.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.396"
$ws.Range("E32").Value = "  -0.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.390"
$ws.Range("E33").Value = "  -2.90%  "
$ws.Range("E34").Value = "  +0.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9933"
$ws.Range("E35").Value = "  -0.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.392"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.751"
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5635"
$ws.Range("E38").Value = "  -5.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01600"
$ws.Range("E39").Value = "  -0.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.869"
$ws.Range("E40").Value = "  -3.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8542"
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("E42").Value = "  -0.32%  "
$ws.Range("D43").Value = "1.032.45"
$ws.Range("E43").Value = "  -6.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "98.59"
$ws.Range("E44").Value = "  -2.03%  "
$ws.Range("D45").Value = "1.793.35"
$ws.Range("E45").Value = "  -1.34%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₈108"
$ws.Range("E46").Value = "  -1.73%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.66"
$ws.Range("E47").Value = "  +0.52%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.000"
$ws.Range("E48").Value = "  -0.21%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.061"
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05143"
$ws.Range("E50").Value = "  -1.84%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4216"
$ws.Range("E51").Value = "  -0.45%  "
